$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 3 (pushes Meagan Waldo and the rest down by one)
$ws.Rows.Item(3).Insert()

# Fill in the new team member's row
$ws.Range("A3").Value = "Johnathon Garcia"
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 0

# The previously-last row (TEAM_MEMBER5) is now row 7; its contribution/attitude
# values are cleared, leaving only the name.
$ws.Range("B7").ClearContents()
$ws.Range("C7").ClearContents()

# Column width adjustments (engine quantizes stored width to the nearest
# 1/6 character unit, so these ColumnWidth inputs are chosen to land on the
# closest achievable stored widths of ~18.14 and 12)
$ws.Columns.Item(1).ColumnWidth = 17.35
$ws.Columns.Item(2).ColumnWidth = 11.15

# Move the selection
$ws.Range("D11").Select() | Out-Null
